$d = $word.ActiveDocument

# --- "Instrukcja.odt" -> "Instrukcja." + "docx" (two runs) ---
$r1 = $d.Content
$r1.Find.Execute("Instrukcja.odt", $false, $false, $false, $false, $false, $true, 1, $false, "")
$full1 = $d.Range($r1.Start, $r1.End)
$tail1 = $d.Range($full1.End - 3, $full1.End)
$tail1.Delete()
$tail1.InsertAfter("docx")
# Toggle a character property on/off so the new text keeps its own run
# instead of being re-merged into the preceding run on save.
$tail1.Bold = 1
$tail1.Bold = 0

# --- "Folder „kod”:" -> "Folder „" + "WykazCmentarza" + "”:" (three runs) ---
$r2 = $d.Content
$r2.Find.Execute("kod”:", $false, $false, $false, $false, $false, $true, 1, $false, "")
$mid2 = $d.Range($r2.Start, $r2.Start + 3)
$mid2.Delete()
$mid2.InsertAfter("WykazCmentarza")
$mid2.Bold = 1
$mid2.Bold = 0
